$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after the existing "2021-Q4"
#    sheet (so the sheet order becomes 2021-Q4, 2022-Q1, 总计).
# ---------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item(1)

$wsQ1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsQ4)
$wsQ1.Name = "2022-Q1"

# NOTE: fetch the "总计" sheet reference *after* the new sheet has been
# inserted (and by name, not position) - Worksheets.Item(2) captured before
# the insert would silently become stale/point at the wrong sheet once the
# collection shifts.
$wsTotal = $wb.Worksheets.Item("总计")

# Re-use the formatting that already exists on the "2021-Q4" sheet: copy the
# header row, then stamp the single data-row format down across the 7 rows
# of data that need to be written.
$wsQ4.Range("A1:H1").Copy($wsQ1.Range("A1:H1"))
$wsQ4.Range("A2:H2").Copy($wsQ1.Range("A2:H2"))
$wsQ4.Range("A2:H2").Copy($wsQ1.Range("A3:H3"))
$wsQ4.Range("A2:H2").Copy($wsQ1.Range("A4:H4"))
$wsQ4.Range("A2:H2").Copy($wsQ1.Range("A5:H5"))
$wsQ4.Range("A2:H2").Copy($wsQ1.Range("A6:H6"))
$wsQ4.Range("A2:H2").Copy($wsQ1.Range("A7:H7"))
$wsQ4.Range("A2:H2").Copy($wsQ1.Range("A8:H8"))

# Columns B-G hold values that look numeric ("506001", "12.84", ...) but must
# stay as text, exactly like the rest of the workbook. Forcing the number
# format to "@" (Text) before the value is assigned stops Excel's automatic
# numeric coercion (and preserves things like leading zeroes).
$wsQ1.Range("B2:G8").NumberFormat = "@"

$fundRows = @(
  @(0, "506001", "万家科创板 2 年定期开放混合型证券投资基金", "12.84", "98.14", "3.95", "0.5072", 5),
  @(1, "005402", "广发资源优选股票A",                         "11.21", "91.29", "4.30", "0.4820", 10),
  @(2, "560003", "益民创新优势混合",                           "4.97",  "77.37", "2.50", "0.1242", 5),
  @(3, "010235", "广发资源优选股票C",                          "2.83",  "91.29", "4.30", "0.1217", 10),
  @(4, "560002", "益民红利成长混合",                           "3.52",  "82.63", "2.44", "0.0859", 8),
  @(5, "005331", "益民优势安享灵活配置混合",                    "1.82",  "41.68", "1.12", "0.0204", 7),
  @(6, "165524", "信诚中证智能家居指数（LOF）",                 "0.40",  "93.89", "1.21", "0.0048", 4)
)

$r = 2
foreach ($row in $fundRows) {
    $wsQ1.Range("A$r").Value = $row[0]
    $wsQ1.Range("B$r").Value = $row[1]
    $wsQ1.Range("C$r").Value = $row[2]
    $wsQ1.Range("D$r").Value = $row[3]
    $wsQ1.Range("E$r").Value = $row[4]
    $wsQ1.Range("F$r").Value = $row[5]
    $wsQ1.Range("G$r").Value = $row[6]
    $wsQ1.Range("H$r").Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: push the existing 2021-Q4 summary row
#    down to row 3 and add a brand-new summary row for 2022-Q1 in row 2.
# ---------------------------------------------------------------------------

# Copy the formatting of the existing summary row down to row 3 first so the
# new row keeps the same look (bold / bordered index cell in column A).
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 3
$wsTotal.Range("D3").Value = 1.03

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 1.35
